# "updated UI for batch operation"
# - Change the "Not Started" default status text to "not_started" everywhere
#   it is used (the whole Swap Status column, C2:C227, shares this string).
# - Reflect the batch-operation UI state: the user scrolled the grid down and
#   multi-selected the full Swap Status column (C2:C227) with C2 as the
#   active cell, ready to apply a bulk edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$statusRange = $ws.Range("C2:C227")
$statusRange.Value = "not_started"

# Scroll the viewport so row 190 is at the top of the window, matching the
# batch-operation view the user left the sheet in.
$win = $excel.ActiveWindow
$win.ScrollRow = 190
$win.ScrollColumn = 1

# Select the whole status column, ready for the batch operation, with C2 as
# the active cell.
$statusRange.Select()
